# Commit: created a way to most likely have different 15 questions until the set runs out
# Adds a new block of 5 quiz questions (rows 47-51) to Sheet1, and fixes a
# duplicate answer-choice in row 42 (C42) by giving it distinct text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix duplicate answer choice B in the "smoking/lung cancer" question (row 42) ---
# C42 used to duplicate B42's text ("random allocation of subjects to smoking is
# unethical."); replace it with a distinct, correct answer choice.
$ws.Range("C42").Value = "observational studies generally cannot rule out confounding."


# --- New question row 47 ---
$ws.Range("A47").Value = "Over the past few years, college enrollments have increased dramatically across the United States. Over this same time period, an educational task force reports that the number of students dropping out of college has significantly increased.
 A college newspaper reports on these results by stating, `"In an effort to meet the demand of more students, it appears that colleges and universities are not providing the same level of support for students to graduate as they were a few years ago.`"
What is wrong with the college newspaper's interpretation
of the results?"
$ws.Range("B47").Value = "Because more people attend college these days, one should expect that more people will drop out of college."
$ws.Range("C47").Value = "The newspaper is assuming that students drop out because of a lack of support."
$ws.Range("D47").Value = "Both A and B are correct."
$ws.Range("E47").Value = "Neither A nor B is correct."
$ws.Range("G47").Value = "a"
$ws.Range("A47").WrapText = $true
$ws.Range("C47").WrapText = $true

# --- New question row 48 ---
$ws.Range("A48").Value = "How can one measure intelligence? One way is to measure the size of the brain. New technology makes it possible to measure the volume of a person's brain in cubic inches without injury. 
What is more, the measurement gives close to the same answer when one repeats it. But how big the brain is has no relation to how smart a person is. As a measure of intelligence, brain volume is"
$ws.Range("B48").Value = "reliable but invalid."
$ws.Range("C48").Value = "valid but not reliable."
$ws.Range("D48").Value = "valid and reliable."
$ws.Range("E48").Value = "not reliable and invalid"
$ws.Range("A48").WrapText = $true

# --- New question row 49 ---
$ws.Range("A49").Value = "A student's research shows that there were more car accidents in 2015 than there were in 1915. He concludes that people were better drivers in 1915 than in 2015. Why is it not valid to use these two numbers to assess driving abilities in these 2 years?"
$ws.Range("B49").Value = "People had more distractions on the road in 2010 than they had in 1910."
$ws.Range("C49").Value = "The numbers were compiled by a student instead of by a professional researcher."
$ws.Range("D49").Value = "The number of cars in the United States increased substantially from 1910 to 2010"
$ws.Range("E49").Value = "One shouldn't compare years that are so far apart."
$ws.Range("G49").Value = "c"
$ws.Range("A49").WrapText = $true

# --- New question row 50 ---
$ws.Range("A50").Value = "In an experiment to study the effect of vibrations on plant growth, the height of a chrysanthemum was measured three times. The reason for making the measurement three times instead of just once was probably to"
$ws.Range("B50").Value = "decrease bias"
$ws.Range("C50").Value = "eliminate confounding."
$ws.Range("D50").Value = "increase reliability"
$ws.Range("E50").Value = "completely eliminate measurement error."
$ws.Range("G50").Value = "c"

# --- New question row 51 ---
$ws.Range("A51").Value = "When repeated measurements each have a systematic error in the same direction, one says that the measurements contain"
$ws.Range("B51").Value = "precision."
$ws.Range("C51").Value = "random error."
$ws.Range("D51").Value = "bias."
$ws.Range("E51").Value = "measurement error"
$ws.Range("G51").Value = "c"

# --- Match the author's final cursor position / selection ---
$null = $ws.Range("C42").Select()
